$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: old "Hammer" weapon entry becomes the new "Speed" item ---
$ws.Range("B5").Value = "Speed"
$ws.Range("C5").Value = "Speed"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 300
$ws.Range("F5:H5").ClearContents()

# --- Row 6: old "Sword" weapon entry becomes the new "Shield" item ---
$ws.Range("B6").Value = "Shield"
$ws.Range("C6").Value = "Shield"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 400
$ws.Range("F6:H6").ClearContents()

# --- Row 7: old "ChainSaw" weapon entry becomes the new "Power" item ---
# (weaponType/damage/rate columns are left untouched on this row)
$ws.Range("B7").Value = "Power"
$ws.Range("C7").Value = "Power"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 500

# --- Row 8: old "DarkSword" weapon entry is cleared out, keeping only the id ---
$ws.Range("B8:H8").ClearContents()

# --- Rows 9-11: new blank placeholder rows, only the itemNum id is filled in ---
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9

# --- Rows 12-15: re-added weapon catalogue (Hammer/Sword/ChainSaw/DarkSword) ---
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Weapon"
$ws.Range("C12").Value = "Hammer"
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 300
$ws.Range("F12").Value = "TwohandSword"
$ws.Range("G12").Value = 20
$ws.Range("H12").Value = 1.3

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Weapon"
$ws.Range("C13").Value = "Sword"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 500
$ws.Range("F13").Value = "Sword"
$ws.Range("G13").Value = 10
$ws.Range("H13").Value = 0.8

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Weapon"
$ws.Range("C14").Value = "ChainSaw"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 800
$ws.Range("F14").Value = "ChainSaw"
$ws.Range("G14").Value = 50
$ws.Range("H14").Value = 2

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Weapon"
$ws.Range("C15").Value = "DarkSword"
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 1000
$ws.Range("F15").Value = "Sword"
$ws.Range("G15").Value = 50
$ws.Range("H15").Value = 0.4

# --- A lone marker value added elsewhere on the sheet ---
$ws.Range("K17").Value = 5

# Restore the selection the author left behind when saving
$ws.Range("B6").Select()
